$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in H1 from "Status" to "Keterangan"
$ws.Range("H1").Value = "Keterangan"

# Update the view: zoom level and active selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$ws.Range("D3").Select() | Out-Null
